$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), styled like the existing
# header cells (bold, centered, bordered) by copying format from H1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the two new columns, rows 2-30.
$iValues = @(6, 6, 1, 6, 6, 7, 10, 8, 7, 5, 7, 8, 8, 9, 9, 3, 7, 6, 8, 7, 7, 6, 5, 5, 5, 5, 5, 5, 3)
$jValues = @(6, 6, 1, 6, 6, 7, 10, 8, 7, 5, 7, 8, 8, 9, 9, 4, 7, 6, 8, 8, 7, 6, 5, 5, 5, 5, 5, 5, 3)

for ($r = 2; $r -le 30; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
